$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 282.72726
$ws.Range("I2").Value = 291
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 291
$ws.Range("L2").Value = 200
$ws.Range("M2").Value = -178
$ws.Range("N2").Value = -426
$ws.Range("H86").Value = 3056
$ws.Range("I86").Value = 2934.1667
$ws.Range("J86").Value = 3299.6667
$ws.Range("K86").Value = 2934.1667
$ws.Range("L86").Value = 3299.6667
$ws.Range("M86").Value = -1811.1667
$ws.Range("N86").Value = -5545.6667
$ws.Range("H89").Value = 3056
$ws.Range("I89").Value = 2934.1667
$ws.Range("J89").Value = 3299.6667
$ws.Range("K89").Value = 14670.8335
$ws.Range("L89").Value = 16498.3335
$ws.Range("M89").Value = -9054.833500000001
$ws.Range("N89").Value = -27730.3335
$ws.Range("H116").Value = 4614517.5
$ws.Range("I116").Value = 5126993.5
$ws.Range("J116").Value = 2233.3333
$ws.Range("K116").Value = 5126993.5
$ws.Range("L116").Value = 2233.3333
$ws.Range("M116").Value = -5123551.5
$ws.Range("N116").Value = -9117.3333
$ws.Range("H125").Value = 7007556
$ws.Range("I125").Value = 675.6923
$ws.Range("J125").Value = 37370704
$ws.Range("K125").Value = 6081.2307
$ws.Range("L125").Value = 336336336
$ws.Range("M125").Value = -3621.2307
$ws.Range("N125").Value = -336341256
$ws.Range("H133").Value = 14976.357
$ws.Range("J133").Value = 14976.357
$ws.Range("L133").Value = 14976.357
$ws.Range("N133").Value = -25096.357

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15283.763
$ws.Range("I32").Value = 2179.3142
$ws.Range("J32").Value = 168169
$ws.Range("K32").Value = 2179.3142
$ws.Range("L32").Value = 168169
$ws.Range("M32").Value = -1892.3142
$ws.Range("N32").Value = -168743
$ws.Range("H61").Value = 3515.64
$ws.Range("I61").Value = 2761.2144
$ws.Range("J61").Value = 4475.8184
$ws.Range("K61").Value = 2761.2144
$ws.Range("L61").Value = 4475.8184
$ws.Range("M61").Value = -2549.2144
$ws.Range("N61").Value = -4899.8184
$ws.Range("H74").Value = 7528.8857
$ws.Range("I74").Value = 1272.7931
$ws.Range("J74").Value = 37766.668
$ws.Range("K74").Value = 1272.7931
$ws.Range("L74").Value = 37766.668
$ws.Range("M74").Value = -398.7931000000001
$ws.Range("N74").Value = -39514.668
$ws.Range("H77").Value = 7528.8857
$ws.Range("I77").Value = 1272.7931
$ws.Range("J77").Value = 37766.668
$ws.Range("K77").Value = 6363.9655
$ws.Range("L77").Value = 188833.34
$ws.Range("M77").Value = -1995.9655
$ws.Range("N77").Value = -197569.34
$ws.Range("H132").Value = 4467.609
$ws.Range("I132").Value = 4389.579
$ws.Range("J132").Value = 4838.25
$ws.Range("K132").Value = 13168.737
$ws.Range("L132").Value = 14514.75
$ws.Range("M132").Value = -10638.737
$ws.Range("N132").Value = -19574.75
$ws.Range("H136").Value = 3515.64
$ws.Range("I136").Value = 2761.2144
$ws.Range("J136").Value = 4475.8184
$ws.Range("K136").Value = 8283.643199999999
$ws.Range("L136").Value = 13427.4552
$ws.Range("M136").Value = -5733.643199999999
$ws.Range("N136").Value = -18527.4552

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 887.4
$ws.Range("I36").Value = 887.4
$ws.Range("K36").Value = 887.4
$ws.Range("M36").Value = -353.4
$ws.Range("H75").Value = 93289.55
$ws.Range("I75").Value = 5952.1665
$ws.Range("K75").Value = 5952.1665
$ws.Range("M75").Value = -5016.1665
$ws.Range("H78").Value = 93289.55
$ws.Range("I78").Value = 5952.1665
$ws.Range("K78").Value = 17856.4995
$ws.Range("M78").Value = -13176.4995
$ws.Range("H134").Value = 3386.0278
$ws.Range("I134").Value = 1912.9546
$ws.Range("J134").Value = 5700.857
$ws.Range("K134").Value = 5738.8638
$ws.Range("L134").Value = 17102.571
$ws.Range("M134").Value = -3203.8638
$ws.Range("N134").Value = -22172.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 109.888885
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 114.833336
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 114.833336
$ws.Range("M7").Value = 13
$ws.Range("N7").Value = -340.833336
$ws.Range("H58").Value = 2176.0417
$ws.Range("I58").Value = 1389.0588
$ws.Range("J58").Value = 4087.2856
$ws.Range("K58").Value = 1389.0588
$ws.Range("L58").Value = 4087.2856
$ws.Range("M58").Value = -1186.0588
$ws.Range("N58").Value = -4493.2856
$ws.Range("H99").Value = 5683237.5
$ws.Range("I99").Value = 6945845.5
$ws.Range("K99").Value = 6945845.5
$ws.Range("M99").Value = -6944347.5
$ws.Range("H126").Value = 5683237.5
$ws.Range("I126").Value = 6945845.5
$ws.Range("K126").Value = 20837536.5
$ws.Range("M126").Value = -20835066.5
$ws.Range("H136").Value = 2176.0417
$ws.Range("I136").Value = 1389.0588
$ws.Range("J136").Value = 4087.2856
$ws.Range("K136").Value = 4167.1764
$ws.Range("L136").Value = 12261.8568
$ws.Range("M136").Value = -1617.1764
$ws.Range("N136").Value = -17361.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2099.5
$ws.Range("I5").Value = 1739.3
$ws.Range("J5").Value = 3000
$ws.Range("K5").Value = 5217.9
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = -5105.9
$ws.Range("N5").Value = -9224
$ws.Range("H80").Value = 1063.25
$ws.Range("J80").Value = 1143.7142
$ws.Range("L80").Value = 3431.1426
$ws.Range("N80").Value = -5303.142599999999
$ws.Range("H83").Value = 1063.25
$ws.Range("J83").Value = 1143.7142
$ws.Range("L83").Value = 10293.4278
$ws.Range("N83").Value = -19653.4278
$ws.Range("H122").Value = 969.2
$ws.Range("I122").Value = 504
$ws.Range("J122").Value = 1020.8889
$ws.Range("K122").Value = 4536
$ws.Range("L122").Value = 9188.000100000001
$ws.Range("M122").Value = -2086
$ws.Range("N122").Value = -14088.0001
$ws.Range("H131").Value = 1466.305
$ws.Range("I131").Value = 365.5
$ws.Range("J131").Value = 1690.9592
$ws.Range("K131").Value = 1096.5
$ws.Range("L131").Value = 5072.8776
$ws.Range("M131").Value = 3943.5
$ws.Range("N131").Value = -15152.8776
$ws.Range("H135").Value = 2099.5
$ws.Range("I135").Value = 1739.3
$ws.Range("J135").Value = 3000
$ws.Range("K135").Value = 15653.7
$ws.Range("L135").Value = 27000
$ws.Range("M135").Value = -13118.7
$ws.Range("N135").Value = -32070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1112439
$ws.Range("I122").Value = 1852805.2
$ws.Range("J122").Value = 1889.5
$ws.Range("K122").Value = 5558415.6
$ws.Range("L122").Value = 5668.5
$ws.Range("M122").Value = -5555965.6
$ws.Range("N122").Value = -10568.5
$ws.Range("H126").Value = 2414
$ws.Range("I126").Value = 1708.3846
$ws.Range("J126").Value = 2953.5881
$ws.Range("K126").Value = 5125.1538
$ws.Range("L126").Value = 8860.764299999999
$ws.Range("M126").Value = -2655.1538
$ws.Range("N126").Value = -13800.7643
$ws.Range("H132").Value = 5130.1816
$ws.Range("I132").Value = 5253
$ws.Range("J132").Value = 4802.6665
$ws.Range("K132").Value = 15759
$ws.Range("L132").Value = 14407.9995
$ws.Range("M132").Value = -13229
$ws.Range("N132").Value = -19467.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2603
$ws.Range("I40").Value = 1107.5294
$ws.Range("J40").Value = 3874.15
$ws.Range("K40").Value = 1107.5294
$ws.Range("L40").Value = 3874.15
$ws.Range("M40").Value = -971.5293999999999
$ws.Range("N40").Value = -4146.15
$ws.Range("H82").Value = 1176.8182
$ws.Range("I82").Value = 998.4
$ws.Range("J82").Value = 1325.5
$ws.Range("K82").Value = 998.4
$ws.Range("L82").Value = 1325.5
$ws.Range("M82").Value = -637.4
$ws.Range("N82").Value = -2047.5
$ws.Range("H85").Value = 1176.8182
$ws.Range("I85").Value = 998.4
$ws.Range("J85").Value = 1325.5
$ws.Range("K85").Value = 998.4
$ws.Range("L85").Value = 1325.5
$ws.Range("M85").Value = 249.6
$ws.Range("N85").Value = -3821.5
$ws.Range("H93").Value = 1384.7273
$ws.Range("I93").Value = 1076.4445
$ws.Range("J93").Value = 2772
$ws.Range("K93").Value = 1076.4445
$ws.Range("L93").Value = 2772
$ws.Range("M93").Value = 171.5554999999999
$ws.Range("N93").Value = -5268
$ws.Range("H136").Value = 5878.3105
$ws.Range("I136").Value = 3747.8
$ws.Range("J136").Value = 10612.777
$ws.Range("K136").Value = 11243.4
$ws.Range("L136").Value = 31838.331
$ws.Range("M136").Value = -8693.400000000001
$ws.Range("N136").Value = -36938.331
$ws.Range("H140").Value = 71738.336
$ws.Range("J140").Value = 71738.336
$ws.Range("L140").Value = 71738.336
$ws.Range("N140").Value = -82098.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 38095.645
$ws.Range("I122").Value = 47197.727
$ws.Range("J122").Value = 4721.3335
$ws.Range("K122").Value = 141593.181
$ws.Range("L122").Value = 14164.0005
$ws.Range("M122").Value = -139143.181
$ws.Range("N122").Value = -19064.0005
$ws.Range("H136").Value = 6430735.5
$ws.Range("I136").Value = 7599291
$ws.Range("J136").Value = 3681.375
$ws.Range("K136").Value = 22797873
$ws.Range("L136").Value = 11044.125
$ws.Range("M136").Value = -22795323
$ws.Range("N136").Value = -16144.125
